$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change
$ws.Range("B3").Value = 0.870628789650563
$ws.Range("C3").Value = 0.8622395829820088
$ws.Range("D3").Value = 0.8459063383793706

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.6944089051338344
$ws.Range("C4").Value = 0.6945929291478993
$ws.Range("D4").Value = 0.4819568061340534

# Row 5: AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.7188741652076042
$ws.Range("C5").Value = 0.670413811078315
$ws.Range("D5").Value = 0.6975499638964551
